# Weekly update for "Fruta, Femacal de La Calera - Tuna" sheet.
# Inserts 3 new rows (a new reporting week, Fecha serial 45021) above the
# existing data block starting at row 195, pushing the previous rows down
# by three positions (195-209 -> 198-212). The dimension / used range grows
# from A1:T209 to A1:T212 automatically once the new rows are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 195; existing rows 195:209
# shift down to 198:212, carrying their formatting (incl. column D's
# date-time number format) with them.
$ws.Rows("195:197").Insert()

# --- New row 195: Tuna, Especial, 16 kilo box ---
$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D195").Value = 45021
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = "Fruta"
$ws.Range("G195").Value = 100107
$ws.Range("H195").Value = "Otros"
$ws.Range("I195").Value = 100107011
$ws.Range("J195").Value = "Tuna"
$ws.Range("K195").Value = "Sin especificar"
$ws.Range("L195").Value = "Especial"
$ws.Range("M195").Value = 56
$ws.Range("N195").Value = 16000
$ws.Range("O195").Value = 16000
$ws.Range("P195").Value = 16000
$ws.Range("Q195").Value = "$/caja 16 kilos"
$ws.Range("R195").Value = "Provincia de Los Andes"
$ws.Range("S195").Value = 1000
$ws.Range("T195").Value = 16

# --- New row 196: Tuna, Primera, 16 kilo box ---
$ws.Range("A196").Value = 3
$ws.Range("B196").Value = "Femacal de La Calera"
$ws.Range("C196").Value = "Coquimbo"
$ws.Range("D196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D196").Value = 45021
$ws.Range("E196").Value = 5
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100107
$ws.Range("H196").Value = "Otros"
$ws.Range("I196").Value = 100107011
$ws.Range("J196").Value = "Tuna"
$ws.Range("K196").Value = "Sin especificar"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 67
$ws.Range("N196").Value = 13000
$ws.Range("O196").Value = 13000
$ws.Range("P196").Value = 13000
$ws.Range("Q196").Value = "$/caja 16 kilos"
$ws.Range("R196").Value = "Provincia de Los Andes"
$ws.Range("S196").Value = 812
$ws.Range("T196").Value = 16

# --- New row 197: Tuna, Segunda, 16 kilo box ---
$ws.Range("A197").Value = 3
$ws.Range("B197").Value = "Femacal de La Calera"
$ws.Range("C197").Value = "Coquimbo"
$ws.Range("D197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D197").Value = 45021
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = "Fruta"
$ws.Range("G197").Value = 100107
$ws.Range("H197").Value = "Otros"
$ws.Range("I197").Value = 100107011
$ws.Range("J197").Value = "Tuna"
$ws.Range("K197").Value = "Sin especificar"
$ws.Range("L197").Value = "Segunda"
$ws.Range("M197").Value = 60
$ws.Range("N197").Value = 10000
$ws.Range("O197").Value = 10000
$ws.Range("P197").Value = 10000
$ws.Range("Q197").Value = "$/caja 16 kilos"
$ws.Range("R197").Value = "Provincia de Los Andes"
$ws.Range("S197").Value = 625
$ws.Range("T197").Value = 16
